$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 4
    3 = 2
    4 = 5
    5 = 3
    6 = 4
    7 = 6
    8 = 2
    9 = 6
    10 = 3
    11 = 4
    12 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
